$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric readings in row 5 to 2 decimal places (custom accuracy).
# Columns A, N, O, S, T, V already only had <=2 decimals and keep their values.
$ws.Range("B5").Value = 4.67
$ws.Range("C5").Value = 3.41
$ws.Range("D5").Value = 0.7
$ws.Range("E5").Value = 10.43
$ws.Range("F5").Value = 7.99
$ws.Range("G5").Value = 3.62
$ws.Range("H5").Value = 17.78
$ws.Range("I5").Value = 5.79
$ws.Range("J5").Value = 2.56
$ws.Range("K5").Value = 3.4
$ws.Range("L5").Value = 4.18
$ws.Range("M5").Value = 4.55
$ws.Range("N5").Value = 1.08
$ws.Range("O5").Value = 3.76
$ws.Range("P5").Value = 5.24
$ws.Range("Q5").Value = 3.41
$ws.Range("R5").Value = 0.64
$ws.Range("S5").Value = 0.25
$ws.Range("T5").Value = 49.92
$ws.Range("U5").Value = 10.67
$ws.Range("V5").Value = 3.47
$ws.Range("W5").Value = 6.94
$ws.Range("X5").Value = 3.61
$ws.Range("Y5").Value = 0.68
$ws.Range("Z5").Value = 8.5
$ws.Range("AA5").Value = 3.07
$ws.Range("AB5").Value = 2.84
$ws.Range("AC5").Value = 3.31
$ws.Range("AD5").Value = 4.39
$ws.Range("AE5").Value = 0.53
$ws.Range("AF5").Value = 16.38
$ws.Range("AG5").Value = 1.83
$ws.Range("AH5").Value = 4.34

# Remove the now-unneeded last data row (row 6), shrinking the sheet
# dimension from A1:AH6 to A1:AH5.
$ws.Rows.Item(6).Delete()
